$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "Finalizado" -> "Cargado" ---
$ws.Range("G1").Value = "Cargado"

# --- New column E (time-of-day durations), reusing the existing h:mm format ---
$ws.Range("E2").NumberFormat = "h:mm"
$ws.Range("E2").Value = 0.020833333333333332
$ws.Range("E3").NumberFormat = "h:mm"
$ws.Range("E3").Value = 0.010416666666666666
$ws.Range("E4").NumberFormat = "h:mm"
$ws.Range("E4").Value = 0.22916666666666666
$ws.Range("E5").NumberFormat = "h:mm"
$ws.Range("E5").Value = 0.020833333333333332

# --- Column G: clear the repeated zeros, merge into a single total cell ---
$ws.Range("G2:G5").HorizontalAlignment = -4108
$ws.Range("G3:G5").ClearContents()
$ws.Range("G2:G5").Merge()
$ws.Range("G2").NumberFormat = "h:mm"
$ws.Range("G2").Value = 0.30208333333333331

# --- Row 6: totals ---
$ws.Range("E6").Formula = "=SUM(E2:E5)"
$ws.Range("F6").Formula = "=SUM(F2:F5)"
$ws.Range("G6").Formula = "=SUM(G2:G5)"
$ws.Range("E6:G6").NumberFormat = "[`$-F400]h:mm:ss\ am/pm"

# --- Column C: widen to fit the longest process name ---
$ws.Columns.Item(3).AutoFit()

# --- Selection moves to F4 ---
[void]$ws.Range("F4").Select()
